$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.435.04"
$ws.Range("E2").Value = "  +0.13%  "
$ws.Range("D3").Value = "1.570.75"
$ws.Range("E3").Value = "  +0.49%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("E5").Value = "  +0.05%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "287.61"
$ws.Range("E6").Value = "  +0.98%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3696"
$ws.Range("E7").Value = "  +2.12%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "47.37"
$ws.Range("E8").Value = "  -2.22%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3317"
$ws.Range("E9").Value = "  -0.58%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.168"
$ws.Range("E10").Value = "  +3.67%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07498"
$ws.Range("E11").Value = "  +1.43%  "
$ws.Range("E12").Value = "  +0.07%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.72"
$ws.Range("E13").Value = "  -0.17%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.930"
$ws.Range("E14").Value = "  -0.03%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.887"
$ws.Range("E15").Value = "  -0.01%  "
$ws.Range("D16").Value = "1.559.85"
$ws.Range("E16").Value = "  -0.34%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001113"
$ws.Range("E17").Value = "  +0.77%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "87.80"
$ws.Range("E18").Value = "  -0.36%  "
$ws.Range("E19").Value = "  +0.46%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.422"
$ws.Range("E20").Value = "  +1.40%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9996"
$ws.Range("E21").Value = "  -0.07%  "
$ws.Range("E22").Value = "  +2.11%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.96"
$ws.Range("E23").Value = "  -0.47%  "
$ws.Range("D24").Value = "22.423.24"
$ws.Range("E24").Value = "  +0.06%  "
$ws.Range("E25").Value = "  -1.81%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.617"
$ws.Range("E26").Value = "  +2.31%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "150.78"
$ws.Range("E27").Value = "  +0.56%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.52"
$ws.Range("E28").Value = "  +0.69%  "
$ws.Range("E29").Value = "  -1.39%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "124.59"
$ws.Range("E30").Value = "  +1.10%  "
$ws.Range("D31").Value = "1.737.06"
$ws.Range("E31").Value = "  -0.14%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.081"
$ws.Range("E32").Value = "  +2.82%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.089"
$ws.Range("E33").Value = "  -0.70%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.973"
$ws.Range("E34").Value = "  -1.05%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.896"
$ws.Range("E35").Value = "  +1.03%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08308"
$ws.Range("E36").Value = "  +0.32%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02432"
$ws.Range("E37").Value = "  +1.18%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06373"
$ws.Range("E38").Value = "  -0.10%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.296"
$ws.Range("E39").Value = "  +0.64%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2208"
$ws.Range("E40").Value = "  -0.13%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.322"
$ws.Range("E41").Value = "  +0.08%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.38"
$ws.Range("E42").Value = "  +2.05%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6216"
$ws.Range("E43").Value = "  +2.26%  "
$ws.Range("E44").Value = "  +0.10%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.97"
$ws.Range("E45").Value = "  +1.53%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6029"
$ws.Range("E46").Value = "  +4.53%  "
$ws.Range("E47").Value = "  +0.56%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.037"
$ws.Range("E48").Value = "  +1.16%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "124.61"
$ws.Range("E49").Value = "  +0.08%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.194"
$ws.Range("E50").Value = "  -1.61%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07189"
$ws.Range("E51").Value = "  -0.23%  "
